$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting existing rows 27-69 down to 28-70.
$ws.Rows(27).Insert()

# Populate the newly inserted row 27 with the new record (same market/product
# metadata as the surrounding rows, new date/quality/volume/price data).
$ws.Cells.Item(27, 1).Value = 5
$ws.Cells.Item(27, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(27, 3).Value = "Maule"
$ws.Cells.Item(27, 4).Value = 44803
$ws.Cells.Item(27, 5).Value = 7
$ws.Cells.Item(27, 6).Value = "Fruta"
$ws.Cells.Item(27, 7).Value = 100107
$ws.Cells.Item(27, 8).Value = "Otros"
$ws.Cells.Item(27, 9).Value = 100107002
$ws.Cells.Item(27, 10).Value = "Chirimoya"
$ws.Cells.Item(27, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(27, 12).Value = "Especial"
$ws.Cells.Item(27, 13).Value = 40
$ws.Cells.Item(27, 14).Value = 30000
$ws.Cells.Item(27, 15).Value = 30000
$ws.Cells.Item(27, 16).Value = 30000
$ws.Cells.Item(27, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(27, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(27, 19).Value = 3000
$ws.Cells.Item(27, 20).Value = 10
